$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New translation rows to append (Key, Polish, English)
$rows = @(
    @("NOINTERNETCONNECTION", "Brak połączenia z Internetem", "No Internet connection"),
    @("INCORRECTMAIL", "Niepoprawny e-mail", "Incorrect e-mail address"),
    @("INVALIDEMAILORPASSWORD", "Niepoprawne dane", "Incorrect data"),
    @("EASYQUIZQUESTION", "Jak nazywa się zaznaczony element?", "What is the name of the selected item?"),
    @("PRESSENTERTOCONFIRM", "Naciśnij ENTER, aby potwierdzić", "Press ENTER to confirm"),
    @("EMAILADDRESSNOTAVAILABLE", "Adres e-mail jest zajęty", "The email address is used")
)

$startRow = 68
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
}

# Update the active selection as in the diff (deselect element after answer in selection quiz)
$ws.Range("B61:B62").Select()
